$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "Voicemail" -> "Voice Mail" wording (rows 5 & 6) ---
$ws.Range("B5").Value = "With Voice Mail Plan"
$ws.Range("B6").Value = "Without Voice Mail Plan"

# --- Replace the "Churn Rate by top_states_by_churn" rows (7-16) with new
#     "Tipping Point" rows, clearing the now-unused churned_count /
#     total_count columns (D & E) for each of them ---

$tippingRows = @(
    @{ Row = 7;  A = "Tipping Point: Total Intl Minutes";     B = "Threshold: 3.9";   C = 14.68 },
    @{ Row = 8;  A = "Tipping Point: Total Intl Charge";      B = "Threshold: 1.05";  C = 14.68 },
    @{ Row = 9;  A = "Tipping Point: Total Day Minutes";      B = "Threshold: 291.2"; C = 75.38 },
    @{ Row = 10; A = "Tipping Point: Total Day Charge";       B = "Threshold: 49.5";  C = 75.38 },
    @{ Row = 11; A = "Tipping Point: Customer Service Calls"; B = "Threshold: 5";     C = 61.39 },
    @{ Row = 12; A = "Tipping Point: International Plan";    B = "Threshold: 1";     C = 42.41 },
    @{ Row = 13; A = "Tipping Point: Total Night Minutes";    B = "Threshold: 104.9"; C = 14.77 },
    @{ Row = 14; A = "Tipping Point: Total Night Charge";     B = "Threshold: 4.72";  C = 14.76 },
    @{ Row = 15; A = "Tipping Point: Total Eve Minutes";      B = "Threshold: 301.0"; C = 29.76 },
    @{ Row = 16; A = "Tipping Point: Total Eve Charge";       B = "Threshold: 25.59"; C = 29.76 }
)

foreach ($r in $tippingRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = ""
    $ws.Cells.Item($r.Row, 5).Value = ""
}

# --- Update the two existing tipping-point rows (17 & 18) with new metrics ---
$ws.Cells.Item(17, 1).Value = "Tipping Point: Account Length"
$ws.Cells.Item(17, 2).Value = "Threshold: 17"
$ws.Cells.Item(17, 3).Value = 14.58

$ws.Cells.Item(18, 1).Value = "Tipping Point: Total Day Calls"
$ws.Cells.Item(18, 2).Value = "Threshold: 141"
$ws.Cells.Item(18, 3).Value = 20

# --- Append new model-performance / risk rows (19-22) ---
$ws.Cells.Item(19, 1).Value = "Model Accuracy"
$ws.Cells.Item(19, 2).Value = "sklearn.LogisticRegression"
$ws.Cells.Item(19, 3).Value = 86.34999999999999
$ws.Cells.Item(19, 4).Value = ""
$ws.Cells.Item(19, 5).Value = ""

$ws.Cells.Item(20, 1).Value = "Model Precision"
$ws.Cells.Item(20, 2).Value = "sklearn.metrics"
$ws.Cells.Item(20, 3).Value = 57.69
$ws.Cells.Item(20, 4).Value = ""
$ws.Cells.Item(20, 5).Value = ""

$ws.Cells.Item(21, 1).Value = "Model Recall"
$ws.Cells.Item(21, 2).Value = "sklearn.metrics"
$ws.Cells.Item(21, 3).Value = 21.74
$ws.Cells.Item(21, 4).Value = ""
$ws.Cells.Item(21, 5).Value = ""

$ws.Cells.Item(22, 1).Value = "High Risk Count"
$ws.Cells.Item(22, 2).Value = "prob > 0.5"
$ws.Cells.Item(22, 3).Value = 834
$ws.Cells.Item(22, 4).Value = ""
$ws.Cells.Item(22, 5).Value = ""
